# daily_report.xlsx - 자동 업데이트: 2025-04-17 00:59:41
# Appends a new week of rows (42-46) to the daily tracking table and
# restores the trailing blank template row (now row 47), mirroring the
# weekly "fill the next empty rows, append one more blank row" pattern
# already used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Carry the existing number formats (date / 2-decimal / integer
#    styles already used by column B, G:I, J:K and N) down across rows
#    42-47 first, so the newly written cells - and the still-blank
#    template row 47 - pick up the same style ids as the rest of the
#    table instead of minting new ones.
# ---------------------------------------------------------------------
$ws.Range("B41").Copy() | Out-Null
$ws.Range("B42:B47").PasteSpecial(-4122) | Out-Null

$ws.Range("G41:I41").Copy() | Out-Null
$ws.Range("G42:I47").PasteSpecial(-4122) | Out-Null

$ws.Range("J41:K41").Copy() | Out-Null
$ws.Range("J42:K47").PasteSpecial(-4122) | Out-Null

$ws.Range("N41").Copy() | Out-Null
$ws.Range("N42:N47").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Raw (non-formula) data for the five new daily rows: 2025-04-12
#    through 2025-04-16.
# ---------------------------------------------------------------------
$ws.Range("A42").Value2 = 41
$ws.Range("A43").Value2 = 42
$ws.Range("A44").Value2 = 43
$ws.Range("A45").Value2 = 44
$ws.Range("A46").Value2 = 45

$ws.Range("B42").Value2 = 45759
$ws.Range("B43").Value2 = 45760
$ws.Range("B44").Value2 = 45761
$ws.Range("B45").Value2 = 45762
$ws.Range("B46").Value2 = 45763

$ws.Range("C42").Value2 = "Sat"
$ws.Range("C43").Value2 = "Sun"
$ws.Range("C44").Value2 = "Mon"
$ws.Range("C45").Value2 = "Tue"
$ws.Range("C46").Value2 = "Wed"

$ws.Range("D42").Value2 = 25029
$ws.Range("D43").Value2 = 25088
$ws.Range("D44").Value2 = 25195
$ws.Range("D45").Value2 = 25835
$ws.Range("D46").Value2 = 25600

$ws.Range("E42").Value2 = 25088
$ws.Range("E43").Value2 = 25195
$ws.Range("E44").Value2 = 25835
$ws.Range("E45").Value2 = 25600
$ws.Range("E46").Value2 = 26102

$ws.Range("L42:L46").Value2 = -10000

$ws.Range("O42").Value2 = 81812
$ws.Range("O43").Value2 = 84763.1
$ws.Range("O44").Value2 = 84024
$ws.Range("O45").Value2 = 85364.7
$ws.Range("O46").Value2 = 85100.2

$ws.Range("P42").Value2 = 84763.1
$ws.Range("P43").Value2 = 84024
$ws.Range("P44").Value2 = 85364.7
$ws.Range("P45").Value2 = 85100.2
$ws.Range("P46").Value2 = 84562.7

# ---------------------------------------------------------------------
# 3) Formulas for the same five rows, written as one fill per column so
#    Excel collapses each into a single shared-formula group (as it
#    already does for every other column block in this sheet).
# ---------------------------------------------------------------------
$ws.Range("F42:F46").Formula = "=E42-D42"
$ws.Range("G42:G46").Formula = "=(E42-`$D`$2)/A42"
$ws.Range("H42:H46").Formula = "=(E42/D42-1)*100"
$ws.Range("I42:I46").Formula = "=(POWER((E42/`$D`$3),1/A42)-1)*100"
$ws.Range("J42:J46").Formula = "=J41*1.013"
$ws.Range("K42:K46").Formula = "=M42-J42"
$ws.Range("M42:M46").Formula = "=L42+E42"
$ws.Range("N42:N46").Formula = "=M42/`$D`$2*100"
$ws.Range("Q42:Q46").Formula = "=P42/`$O`$2*100"

# Row 47 stays the blank "next entry" template row (styles only, copied
# in step 1 above) - matching the shape rows 42-44 had before this edit.

# ---------------------------------------------------------------------
# 4) View state: move the live selection and refresh the zoom level.
# ---------------------------------------------------------------------
$ws.Range("N33").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100

Write-Host "daily_report rows 42-47 updated"
